$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.717.49"
$ws.Range("E2").Value = "  +2.38%  "

$ws.Range("D3").Value = "3.371.17"
$ws.Range("E3").Value = "  +2.25%  "

$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "557.85"
$ws.Range("E5").Value = "  +2.55%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "174.74"
$ws.Range("E6").Value = "  +2.09%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.622"
$ws.Range("E7").Value = "  +2.18%  "

$ws.Range("D8").Value = "3.362.82"
$ws.Range("E8").Value = "  +2.26%  "

$ws.Range("E9").Value = "  +0.06%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.165"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.630"
$ws.Range("E11").Value = "  +3.57%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.89"
$ws.Range("E12").Value = "  +2.93%  "

$ws.Range("E13").Value = "  +4.33%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.07"
$ws.Range("E14").Value = "  +2.73%  "

$ws.Range("D15").Value = "3.921.08"
$ws.Range("E15").Value = "  +2.71%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.23"
$ws.Range("E16").Value = "  +1.60%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.119"
$ws.Range("E17").Value = "  +2.46%  "

$ws.Range("D18").Value = "3.389.07"
$ws.Range("E18").Value = "  +2.87%  "

$ws.Range("D19").Value = "64.655.23"
$ws.Range("E19").Value = "  +2.59%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.73"
$ws.Range("E20").Value = "  +1.17%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.989"
$ws.Range("E21").Value = "  +2.73%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "459.26"
$ws.Range("E22").Value = "  +8.58%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.85"
$ws.Range("E23").Value = "  +9.07%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.11"
$ws.Range("E24").Value = "  +2.32%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.17"
$ws.Range("E25").Value = "  +4.25%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.51"
$ws.Range("E26").Value = "  +2.15%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.93"
$ws.Range("E27").Value = "  +8.18%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.77"
$ws.Range("E28").Value = "  +2.17%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.71"
$ws.Range("E29").Value = "  +1.66%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.60"
$ws.Range("E30").Value = "  +5.65%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.68"
$ws.Range("E31").Value = "  +5.69%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.42"
$ws.Range("E32").Value = "  +1.34%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "569.93"
$ws.Range("E33").Value = "  -0.21%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "61.10"
$ws.Range("E34").Value = "  +5.65%  "

$ws.Range("E35").Value = "  +2.15%  "

$ws.Range("E36").Value = "  +0.00%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.64"
$ws.Range("E37").Value = "  +6.95%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.139"
$ws.Range("E38").Value = "  -3.81%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.37"
$ws.Range("E39").Value = "  +1.40%  "

$ws.Range("D40").Value = "0.0₃0738"
$ws.Range("E40").Value = "  +0.39%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.368"
$ws.Range("E41").Value = "  +1.84%  "

$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.17%  "

$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "3.070.93"
$ws.Range("E43").Value = "  -1.08%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.82"
$ws.Range("E44").Value = "  +2.65%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0414"
$ws.Range("E45").Value = "  +3.85%  "

$ws.Range("E46").Value = "  +4.61%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.44"
$ws.Range("E47").Value = "  +1.36%  "

$ws.Range("E48").Value = "  -2.66%  "

$ws.Range("E49").Value = "  +1.02%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "139.51"
$ws.Range("E50").Value = "  +5.72%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.19"
$ws.Range("E51").Value = "  +2.31%  "
